$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.960.45'
$ws.Range('D3').Value = '2.418.02'
$ws.Range('E3').Value = '  +2.04%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.68%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.85'
$ws.Range('E5').Value = '  +2.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.43'
$ws.Range('E6').Value = '  +5.48%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.541'
$ws.Range('E8').Value = '  +2.38%  '
$ws.Range('D9').Value = '2.450.36'
$ws.Range('E9').Value = '  +3.58%  '
$ws.Range('E10').Value = '  +6.21%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('E12').Value = '  +3.11%  '
$ws.Range('E13').Value = '  +5.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.21'
$ws.Range('E14').Value = '  +6.72%  '
$ws.Range('E15').Value = '  +7.92%  '
$ws.Range('D16').Value = '2.978.91'
$ws.Range('E16').Value = '  +6.43%  '
$ws.Range('D17').Value = '62.820.58'
$ws.Range('E17').Value = '  +5.70%  '
$ws.Range('D18').Value = '2.466.69'
$ws.Range('E18').Value = '  +4.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.90'
$ws.Range('E19').Value = '  -2.50%  '
$ws.Range('E20').Value = '  +4.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '328.04'
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('E22').Value = '  +2.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.04'
$ws.Range('E23').Value = '  +12.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.71'
$ws.Range('E25').Value = '  +2.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '621.55'
$ws.Range('E26').Value = '  +12.27%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.10'
$ws.Range('E27').Value = '  +10.23%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.40'
$ws.Range('E28').Value = '  +3.78%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0981'
$ws.Range('E29').Value = '  +7.04%  '
$ws.Range('B30').Value = 'WrappedeETH'
$ws.Range('C30').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D30').Value = '2.565.38'
$ws.Range('E30').Value = '  +3.24%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.15'
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('E32').Value = '  +8.63%  '
$ws.Range('E33').Value = '  +6.03%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +4.02%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.48'
$ws.Range('E35').Value = '  +4.98%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.995'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.75'
$ws.Range('E37').Value = '  +5.00%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.373'
$ws.Range('E38').Value = '  +2.15%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '152.76'
$ws.Range('E39').Value = '  +1.10%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.39'
$ws.Range('E40').Value = '  +8.05%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.67'
$ws.Range('E41').Value = '  +3.22%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.75'
$ws.Range('E42').Value = '  +15.24%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.76'
$ws.Range('E43').Value = '  +7.28%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '0.0₆0284'
$ws.Range('E45').Value = '  -5.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '144.89'
$ws.Range('E46').Value = '  +5.01%  '
$ws.Range('E47').Value = '  +2.57%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.25'
$ws.Range('E48').Value = '  +6.12%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.602'
$ws.Range('E49').Value = '  +3.05%  '
$ws.Range('E50').Value = '  +3.40%  '
$ws.Range('E51').Value = '  +2.86%  '
